$d = $word.ActiveDocument

# 1) ", portador do CNPJ:"  ->  ", inscrito no CNPJ:"
#    (plain text swap inside the existing non-bold run; formatting is
#    untouched because the replacement stays within that run's range)
$r1 = $d.Content
$found1 = $r1.Find.Execute(", portador do CNPJ:", $true, $false, $false, $false, $false, $true, 1, $false, ", inscrito no CNPJ:", 2)
if (-not $found1) { throw "could not find ', portador do CNPJ:'" }

# 2) " #CNPJ"  ->  " #CNPJ " (append a trailing space; stays in the bold
#    "#CNPJ" run so the text keeps its bold/bCs formatting)
$r2 = $d.Content
$found2 = $r2.Find.Execute(" #CNPJ", $true, $false, $false, $false, $false, $true, 1, $false, " #CNPJ ", 2)
if (-not $found2) { throw "could not find ' #CNPJ'" }

# 3) Insert the new address/cep placeholders (non-bold) right before
#    ", representado por". Prepending the new text via Find/Replace on
#    that run keeps its (non-bold) formatting -- inserting via
#    Collapse+InsertAfter on the previous (bold) run would instead have
#    inherited the bold formatting from the left.
$r3 = $d.Content
$found3 = $r3.Find.Execute(", representado por", $true, $false, $false, $false, $false, $true, 1, $false, "e com sede na #END_EMPRESA cep: #CP_EMPRESA, representado por", 2)
if (-not $found3) { throw "could not find ', representado por'" }

Write-Host "found1=$found1 found2=$found2 found3=$found3"
